# Weekly update: insert the new week's price record for Jengibre
# (Vega Central Mapocho de Santiago) as a new row 15, pushing all the
# existing rows (old rows 15-121) down by one (to 16-122).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 15; existing row 15 (and all rows
# below it) shift down by one row.
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with this week's data point.
$ws.Range("A15").Value = 9
$ws.Range("B15").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C15").Value = "Metropolitana"
$ws.Range("D15").Value = 44970
$ws.Range("E15").Value = 13
$ws.Range("F15").Value = 100114007
$ws.Range("G15").Value = "Jengibre"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 520
$ws.Range("K15").Value = 19000
$ws.Range("L15").Value = 20000
$ws.Range("M15").Value = 19500
$ws.Range("N15").Value = "$/caja 13 kilos"
$ws.Range("O15").Value = "Perú"
$ws.Range("P15").Value = 1500
$ws.Range("Q15").Value = 13
$ws.Range("R15").Value = "Hortaliza"
